$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 77, shifting existing rows 77-78 down to 78-79.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with this week's price data.
$ws.Cells.Item(77, 1).Value = 8
$ws.Cells.Item(77, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(77, 3).Value = "Coquimbo"
$ws.Cells.Item(77, 4).Value = 44448
$ws.Cells.Item(77, 5).Value = 4
$ws.Cells.Item(77, 6).Value = 100112040
$ws.Cells.Item(77, 7).Value = "Cilantro"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 3400
$ws.Cells.Item(77, 11).Value = 2000
$ws.Cells.Item(77, 12).Value = 2500
$ws.Cells.Item(77, 13).Value = 2250
$ws.Cells.Item(77, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(77, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(77, 16).Value = 1500
$ws.Cells.Item(77, 17).Value = 1.5
$ws.Cells.Item(77, 18).Value = "Hortaliza"
